$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total hours (C1) changes from 30 to 24
$ws.Range("C1").Value = 24

# Update the coverage formulas in column E (rows 3-16) to divide by $C$1
# instead of the hardcoded 30
for ($r = 3; $r -le 16; $r++) {
    $ws.Range("E$r").Formula = "=SUM(`$C`$3:C$r)/`$C`$1"
}

# New lab ("AutoML frameworks") added at row 14, pushing the "Challenge"
# row down to row 15
$ws.Range("B14").Value = "AutoML frameworks"
$ws.Range("B15").Value = "Challenge"

# Row 16 no longer carries a "What" label (it moved to row 14)
$ws.Range("B16").Clear()

# Row 17 is emptied out entirely: A17/E17 keep their formatting but lose
# their content, while B17/C17 are fully cleared (formatting included)
$ws.Range("A17").ClearContents()
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("E17").ClearContents()

# Update the active selection
$ws.Range("B13").Select()
